$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in row 9 (Exp 13) with its parameters
$ws.Range("A9").Value = "Exp 13"
$ws.Range("B9").Value = 0.75
$ws.Range("C9").Value = 1
$ws.Range("F9").Value = "Exp 13.png"

# Move the active selection to E16, matching the saved view state
$ws.Range("E16").Select()
